# Update workbook data to match the "output generated at 456a3b4" refresh.
# Numeric "想去人数" (want-to-go count) values increased for many rows, and
# two Cover image URLs were refreshed, across sheets 展览 and 全部类型
# (plus single isolated bumps on 演出 and 本地生活).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
foreach ($chg in @(
    @{ Cell = "F2";  Value = 724 }
    @{ Cell = "F3";  Value = 64 }
    @{ Cell = "F4";  Value = 2027 }
    @{ Cell = "F5";  Value = 5989 }
    @{ Cell = "I5";  Value = "//i2.hdslb.com/bfs/openplatform/202404/fa41FZWy1714473760924.jpeg" }
    @{ Cell = "F6";  Value = 1689 }
    @{ Cell = "F8";  Value = 3373 }
    @{ Cell = "F11"; Value = 1411 }
    @{ Cell = "F12"; Value = 4696 }
    @{ Cell = "F13"; Value = 1768 }
    @{ Cell = "I14"; Value = "//i0.hdslb.com/bfs/openplatform/202404/WR7lMMzi1714474414048.jpeg" }
    @{ Cell = "F17"; Value = 213 }
    @{ Cell = "F19"; Value = 1050 }
    @{ Cell = "F27"; Value = 5 }
    @{ Cell = "F30"; Value = 116 }
    @{ Cell = "F31"; Value = 227 }
    @{ Cell = "F34"; Value = 1815 }
    @{ Cell = "F35"; Value = 2292 }
    @{ Cell = "F36"; Value = 1076 }
    @{ Cell = "F38"; Value = 8 }
    @{ Cell = "F39"; Value = 292 }
    @{ Cell = "F40"; Value = 24 }
    @{ Cell = "F41"; Value = 655 }
    @{ Cell = "F42"; Value = 435 }
    @{ Cell = "F43"; Value = 55 }
    @{ Cell = "F44"; Value = 688 }
    @{ Cell = "F45"; Value = 47 }
    @{ Cell = "F46"; Value = 463 }
    @{ Cell = "F47"; Value = 456 }
    @{ Cell = "F49"; Value = 154 }
)) {
    $ws.Range($chg.Cell).Value = $chg.Value
}

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
foreach ($chg in @(
    @{ Cell = "F26"; Value = 31 }
)) {
    $ws.Range($chg.Cell).Value = $chg.Value
}

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
foreach ($chg in @(
    @{ Cell = "F2"; Value = 811 }
)) {
    $ws.Range($chg.Cell).Value = $chg.Value
}

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
foreach ($chg in @(
    @{ Cell = "F2";  Value = 726 }
    @{ Cell = "F3";  Value = 64 }
    @{ Cell = "F4";  Value = 2027 }
    @{ Cell = "F5";  Value = 5989 }
    @{ Cell = "I5";  Value = "//i2.hdslb.com/bfs/openplatform/202404/fa41FZWy1714473760924.jpeg" }
    @{ Cell = "F6";  Value = 1689 }
    @{ Cell = "F9";  Value = 3373 }
    @{ Cell = "F11"; Value = 1411 }
    @{ Cell = "F12"; Value = 4696 }
    @{ Cell = "F13"; Value = 1768 }
    @{ Cell = "I14"; Value = "//i0.hdslb.com/bfs/openplatform/202404/WR7lMMzi1714474414048.jpeg" }
    @{ Cell = "F21"; Value = 213 }
    @{ Cell = "F24"; Value = 1050 }
    @{ Cell = "F33"; Value = 116 }
    @{ Cell = "F34"; Value = 227 }
    @{ Cell = "F37"; Value = 1815 }
    @{ Cell = "F38"; Value = 2292 }
    @{ Cell = "F43"; Value = 292 }
    @{ Cell = "F46"; Value = 688 }
    @{ Cell = "F47"; Value = 463 }
    @{ Cell = "F48"; Value = 456 }
)) {
    $ws.Range($chg.Cell).Value = $chg.Value
}

$wb.Save()
